# edit.ps1 - apply the Zadanie2.docx changes:
#  1) Split the run "przy różnych wczytywanych z wejścia punktach startowych ("
#     into "przy" + " różnych wczytywanych z wejścia punktach startowych ("
#     (mirrors the gramStart/gramEnd proof-error split introduced upstream).
#  2) Append a new sentence after the manual line break (<w:cr/>) in the
#     "Metoda obliczeniowa: ..." paragraph.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: split "przy różnych wczytywanych z wejścia punktach startowych ("
# ---------------------------------------------------------------------
$target = "przy różnych wczytywanych z wejścia punktach startowych ("

$findRange = $d.Content
$findRange.Find.ClearFormatting()
$found = $findRange.Find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $sentenceStart = $findRange.Start

    # Re-type the leading word "przy" on its own so the run splits exactly
    # the way Word splits it when it flags the word for grammar review.
    $wordRange = $d.Range($sentenceStart, $sentenceStart + 4)
    $wordRange.Text = ""
    $reinsert = $d.Range($sentenceStart, $sentenceStart)
    $reinsert.InsertBefore("przy")
}

# ---------------------------------------------------------------------
# Change 2: add the new requirements sentence after the manual line break
# ---------------------------------------------------------------------
$crRange = $d.Content
$crRange.Find.ClearFormatting()
$foundCr = $crRange.Find.Execute("ów równań", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($foundCr) {
    $para = $crRange.Paragraphs(1)
    $insertPos = $para.Range.End - 1
    $insertPoint = $d.Range($insertPos, $insertPos)
    $insertPoint.InsertAfter("wprowadzenie przez użytkownika: dokładność, ilość krokow,")
}
